$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the to-do text to also mention player disconnect notifications.
$ws.Range("A19").Value = "Добавить оповещение о присоединении и выходе игрока"

# Move the active selection to match the author's saved view state
# (selection moved from B17 to B23 as rows were reviewed/scrolled through).
$ws.Range("B23").Select()
